$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 140.125
$ws.Range("I9").Value = 130.16667
$ws.Range("J9").Value = 170
$ws.Range("K9").Value = 130.16667
$ws.Range("L9").Value = 170
$ws.Range("M9").Value = 38.83332999999999
$ws.Range("N9").Value = -508
$ws.Range("H43").Value = 7938651
$ws.Range("I43").Value = 1000
$ws.Range("J43").Value = 11113711
$ws.Range("K43").Value = 1000
$ws.Range("L43").Value = 11113711
$ws.Range("M43").Value = -931
$ws.Range("N43").Value = -11113849
$ws.Range("H87").Value = 27980.727
$ws.Range("J87").Value = 27980.727
$ws.Range("L87").Value = 27980.727
$ws.Range("N87").Value = -30476.727
$ws.Range("H90").Value = 27980.727
$ws.Range("J90").Value = 27980.727
$ws.Range("L90").Value = 83942.181
$ws.Range("N90").Value = -96422.181
$ws.Range("H129").Value = 796.1795
$ws.Range("J129").Value = 868.96875
$ws.Range("L129").Value = 2606.90625
$ws.Range("N129").Value = -12606.90625
$ws.Range("H138").Value = 2298.247
$ws.Range("I138").Value = 3196.4285
$ws.Range("J138").Value = 2221.5732
$ws.Range("K138").Value = 9589.2855
$ws.Range("L138").Value = 6664.719599999999
$ws.Range("M138").Value = -4449.2855
$ws.Range("N138").Value = -16944.7196

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H8").Value = 5000400
$ws.Range("I8").Value = 5000400
$ws.Range("J8").Value = 0
$ws.Range("K8").Value = 5000400
$ws.Range("L8").Value = 0
$ws.Range("N8").ClearContents()
$ws.Range("M8").Value = -5000256
$ws.Range("H32").Value = 7148.9355
$ws.Range("I32").Value = 6326.1743
$ws.Range("J32").Value = 17257.143
$ws.Range("K32").Value = 6326.1743
$ws.Range("L32").Value = 17257.143
$ws.Range("M32").Value = -6039.1743
$ws.Range("N32").Value = -17831.143
$ws.Range("H61").Value = 27028108
$ws.Range("I61").Value = 30303920
$ws.Range("J61").Value = 2653.5
$ws.Range("K61").Value = 30303920
$ws.Range("L61").Value = 2653.5
$ws.Range("M61").Value = -30303708
$ws.Range("N61").Value = -3077.5
$ws.Range("H74").Value = 3069.1538
$ws.Range("I74").Value = 2010.625
$ws.Range("J74").Value = 4762.8
$ws.Range("K74").Value = 2010.625
$ws.Range("L74").Value = 4762.8
$ws.Range("M74").Value = -1136.625
$ws.Range("N74").Value = -6510.8
$ws.Range("H77").Value = 3069.1538
$ws.Range("I77").Value = 2010.625
$ws.Range("J77").Value = 4762.8
$ws.Range("K77").Value = 10053.125
$ws.Range("L77").Value = 23814
$ws.Range("M77").Value = -5685.125
$ws.Range("N77").Value = -32550
$ws.Range("H88").Value = 2101
$ws.Range("I88").Value = 1099.6666
$ws.Range("J88").Value = 3102.3333
$ws.Range("K88").Value = 1099.6666
$ws.Range("L88").Value = 3102.3333
$ws.Range("M88").Value = -693.6666
$ws.Range("N88").Value = -3914.3333
$ws.Range("H91").Value = 2101
$ws.Range("I91").Value = 1099.6666
$ws.Range("J91").Value = 3102.3333
$ws.Range("K91").Value = 1099.6666
$ws.Range("L91").Value = 3102.3333
$ws.Range("M91").Value = 304.3334
$ws.Range("N91").Value = -5910.3333
$ws.Range("H97").Value = 613.9
$ws.Range("I97").Value = 613.9
$ws.Range("K97").Value = 613.9
$ws.Range("M97").Value = -117.9
$ws.Range("H110").Value = 1131.4445
$ws.Range("I110").Value = 381.42856
$ws.Range("K110").Value = 381.42856
$ws.Range("M110").Value = 1663.57144
$ws.Range("H113").Value = 22199
$ws.Range("J113").Value = 22199
$ws.Range("L113").Value = 22199
$ws.Range("N113").Value = -30877
$ws.Range("H122").Value = 1923.1666
$ws.Range("J122").Value = 2209.75
$ws.Range("L122").Value = 6629.25
$ws.Range("N122").Value = -11529.25
$ws.Range("H132").Value = 2254.651
$ws.Range("I132").Value = 1592.6279
$ws.Range("J132").Value = 3678
$ws.Range("K132").Value = 4777.8837
$ws.Range("L132").Value = 11034
$ws.Range("M132").Value = -2247.8837
$ws.Range("N132").Value = -16094
$ws.Range("H136").Value = 27028108
$ws.Range("I136").Value = 30303920
$ws.Range("J136").Value = 2653.5
$ws.Range("K136").Value = 90911760
$ws.Range("L136").Value = 7960.5
$ws.Range("M136").Value = -90909210
$ws.Range("N136").Value = -13060.5

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1753.35
$ws.Range("I20").Value = 1858.7333
$ws.Range("K20").Value = 1858.7333
$ws.Range("M20").Value = -1611.7333
$ws.Range("H134").Value = 864.3
$ws.Range("I134").Value = 825.13794
$ws.Range("J134").Value = 2000
$ws.Range("K134").Value = 2475.41382
$ws.Range("L134").Value = 6000
$ws.Range("M134").Value = 59.58618000000024
$ws.Range("N134").Value = -11070

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1173.2192
$ws.Range("I31").Value = 1053.5156
$ws.Range("K31").Value = 1053.5156
$ws.Range("M31").Value = -758.5155999999999
$ws.Range("H34").Value = 1173.2192
$ws.Range("I34").Value = 1053.5156
$ws.Range("K34").Value = 1053.5156
$ws.Range("M34").Value = -851.5155999999999
$ws.Range("H107").Value = 507.23077
$ws.Range("I107").Value = 443.5625
$ws.Range("K107").Value = 443.5625
$ws.Range("M107").Value = 1476.4375
$ws.Range("H132").Value = 1725.1724
$ws.Range("I132").Value = 1221.6
$ws.Range("K132").Value = 3664.8
$ws.Range("M132").Value = -1134.8

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 4805.933
$ws.Range("J107").Value = 6087.826
$ws.Range("L107").Value = 18263.478
$ws.Range("N107").Value = -22103.478
$ws.Range("H131").Value = 23846014
$ws.Range("J131").Value = 49264.324
$ws.Range("L131").Value = 147792.972
$ws.Range("N131").Value = -157872.972
$ws.Range("H138").Value = 2935.5144
$ws.Range("I138").Value = 3003.7693
$ws.Range("J138").Value = 2895.182
$ws.Range("K138").Value = 9011.3079
$ws.Range("L138").Value = 8685.545999999998
$ws.Range("M138").Value = -3871.3079
$ws.Range("N138").Value = -18965.546

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2649.6667
$ws.Range("I132").Value = 2498.5386
$ws.Range("J132").Value = 3042.6
$ws.Range("K132").Value = 7495.6158
$ws.Range("L132").Value = 9127.799999999999
$ws.Range("M132").Value = -4965.6158
$ws.Range("N132").Value = -14187.8
$ws.Range("H135").Value = 36938.89
$ws.Range("J135").Value = 35306.25
$ws.Range("L135").Value = 35306.25
$ws.Range("N135").Value = -45446.25

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2810.6667
$ws.Range("J7").Value = 3136
$ws.Range("L7").Value = 3136
$ws.Range("N7").Value = -3360
$ws.Range("H22").Value = 727.5714
$ws.Range("I22").Value = 848.5
$ws.Range("J22").Value = 636.875
$ws.Range("K22").Value = 848.5
$ws.Range("L22").Value = 636.875
$ws.Range("M22").Value = -553.5
$ws.Range("N22").Value = -1226.875
$ws.Range("H27").Value = 727.5714
$ws.Range("I27").Value = 848.5
$ws.Range("J27").Value = 636.875
$ws.Range("K27").Value = 848.5
$ws.Range("L27").Value = 636.875
$ws.Range("M27").Value = -741.5
$ws.Range("N27").Value = -850.875
$ws.Range("H31").Value = 4113.857
$ws.Range("I31").Value = 0
$ws.Range("J31").Value = 4113.857
$ws.Range("K31").Value = 0
$ws.Range("L31").Value = 4113.857
$ws.Range("M31").ClearContents()
$ws.Range("N31").Value = -4609.857
$ws.Range("H46").Value = 5700
$ws.Range("I46").Value = 1000
$ws.Range("J46").Value = 6875
$ws.Range("K46").Value = 1000
$ws.Range("L46").Value = 6875
$ws.Range("M46").Value = -812
$ws.Range("N46").Value = -7251
$ws.Range("H55").Value = 274.90475
$ws.Range("I55").Value = 208.78572
$ws.Range("K55").Value = 208.78572
$ws.Range("M55").Value = -35.78572
$ws.Range("H126").Value = 2810.6667
$ws.Range("J126").Value = 3136
$ws.Range("L126").Value = 9408
$ws.Range("N126").Value = -14348
$ws.Range("H132").Value = 2842.1538
$ws.Range("I132").Value = 2881.6365
$ws.Range("K132").Value = 8644.9095
$ws.Range("M132").Value = -6114.9095

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H123").Value = 30000
$ws.Range("J123").Value = 30000
$ws.Range("L123").Value = 30000
$ws.Range("N123").Value = -39800
$ws.Range("H137").Value = 35040.832
$ws.Range("J137").Value = 35040.832
$ws.Range("L137").Value = 35040.832
$ws.Range("N137").Value = -45240.832
